$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.3740083992081651
$ws.Range("D2").Value = 0.02724418757612312
$ws.Range("E2").Value = 0.1401897643330834
$ws.Range("F2").Value = 6.557428756354085
$ws.Range("G2").Value = 6.019396530403014
$ws.Range("H2").Value = 3.997281369321115
$ws.Range("J2").Value = 0.3938509743047263
$ws.Range("L2").Value = 0.1118087815777127
$ws.Range("N2").Value = 2.873581687696401
# Row 3
$ws.Range("C3").Value = 0.3726670103702929
$ws.Range("D3").Value = 0.02489096839943272
$ws.Range("E3").Value = 0.140626191770691
$ws.Range("F3").Value = 6.43494821266134
$ws.Range("G3").Value = 5.867864785000791
$ws.Range("H3").Value = 3.938905539796167
$ws.Range("J3").Value = 0.3915715013672525
$ws.Range("L3").Value = 0.1124519255948826
$ws.Range("N3").Value = 2.562605684679681
# Row 4
$ws.Range("C4").Value = 0.3720460413837117
$ws.Range("D4").Value = 0.02343810531193213
$ws.Range("E4").Value = 0.1409525226156241
$ws.Range("F4").Value = 6.364209360459199
$ws.Range("G4").Value = 5.779160603376454
$ws.Range("H4").Value = 3.905675968506898
$ws.Range("J4").Value = 0.3903992367864859
$ws.Range("L4").Value = 0.1128767092591527
$ws.Range("N4").Value = 2.371325805375818
# Row 5
$ws.Range("C5").Value = 0.3718436792811559
$ws.Range("D5").Value = 0.02284387653259046
$ws.Range("E5").Value = 0.1411001577117581
$ws.Range("F5").Value = 6.336493176804538
$ws.Range("G5").Value = 5.744089390514944
$ws.Range("H5").Value = 3.89278550608401
$ws.Range("J5").Value = 0.3899783203635465
$ws.Range("L5").Value = 0.1130573438889471
$ws.Range("N5").Value = 2.293303068605894
# Row 6
$ws.Range("C6").Value = 0.3718131309073129
$ws.Range("D6").Value = 0.02274506901950701
$ws.Range("E6").Value = 0.1411255567447682
$ws.Range("F6").Value = 6.331957644482998
$ws.Range("G6").Value = 5.738330452565265
$ws.Range("H6").Value = 3.890684183614837
$ws.Range("J6").Value = 0.3899118467761795
$ws.Range("L6").Value = 0.1130877934985755
$ws.Range("N6").Value = 2.2803432614038
# Row 7
$ws.Range("C7").Value = 0.37204310737323
$ws.Range("D7").Value = 0.0234301003647488
$ws.Range("E7").Value = 0.1409544543751622
$ws.Range("F7").Value = 6.363831089935104
$ws.Range("G7").Value = 5.778683280534892
$ws.Range("H7").Value = 3.905499496123923
$ws.Range("J7").Value = 0.3903933307030911
$ws.Range("L7").Value = 0.1128791148437678
$ws.Range("N7").Value = 2.370273851392596
# Row 8
$ws.Range("C8").Value = 0.3735036620072663
$ws.Range("D8").Value = 0.02643434981806081
$ws.Range("E8").Value = 0.1403281167795765
$ws.Range("F8").Value = 6.514263321210478
$ws.Range("G8").Value = 5.96623961162669
$ws.Range("H8").Value = 3.976607092242034
$ws.Range("J8").Value = 0.393017596131898
$ws.Range("L8").Value = 0.1120243431398293
$ws.Range("N8").Value = 2.766433886209882
# Row 9
$ws.Range("C9").Value = 0.3779889155190688
$ws.Range("D9").Value = 0.03227036071609746
$ws.Range("E9").Value = 0.1395641598401873
$ws.Range("F9").Value = 6.845272920612103
$ws.Range("G9").Value = 6.369124433330228
$ws.Range("H9").Value = 4.137089045987068
$ws.Range("J9").Value = 0.3999858260429789
$ws.Range("L9").Value = 0.1105846253203246
$ws.Range("N9").Value = 3.540180268007646
# Row 10
$ws.Range("C10").Value = 0.3822922488742222
$ws.Range("D10").Value = 0.03653591584334492
$ws.Range("E10").Value = 0.1392878426893951
$ws.Range("F10").Value = 7.111286806379837
$ws.Range("G10").Value = 6.687509928407451
$ws.Range("H10").Value = 4.268271190672465
$ws.Range("J10").Value = 0.4062432343927043
$ws.Range("L10").Value = 0.1096701222095007
$ws.Range("N10").Value = 4.10623028343673
# Row 11
$ws.Range("C11").Value = 0.3844733748431963
$ws.Range("D11").Value = 0.03847420542636826
$ws.Range("E11").Value = 0.1392244798864475
$ws.Range("F11").Value = 7.237455936493802
$ws.Range("G11").Value = 6.837438603795249
$ws.Range("H11").Value = 4.330934473592549
$ws.Range("J11").Value = 0.4093431699043464
$ws.Range("L11").Value = 0.1092850110906998
$ws.Range("N11").Value = 4.363110593465422
# Row 12
$ws.Range("C12").Value = 0.3853318263369374
$ws.Range("D12").Value = 0.03920808962710964
$ws.Range("E12").Value = 0.1392094891314848
$ws.Range("F12").Value = 7.285990733471635
$ws.Range("G12").Value = 6.894963572998677
$ws.Range("H12").Value = 4.355101455187651
$ws.Range("J12").Value = 0.4105539854505622
$ws.Range("L12").Value = 0.1091436088602382
$ws.Range("N12").Value = 4.460285735713398
# Row 13
$ws.Range("C13").Value = 0.3851454917516719
$ws.Range("D13").Value = 0.03905003587206579
$ws.Range("E13").Value = 0.139212316574941
$ws.Range("F13").Value = 7.275503962648315
$ws.Range("G13").Value = 6.882540885201593
$ws.Range("H13").Value = 4.34987706190077
$ws.Range("J13").Value = 0.4102915637443232
$ws.Range("L13").Value = 0.1091738654864134
$ws.Range("N13").Value = 4.439361943450422
# Row 14
$ws.Range("C14").Value = 0.3845433465452288
$ws.Range("D14").Value = 0.03853458342437222
$ws.Range("E14").Value = 0.1392230659115725
$ws.Range("F14").Value = 7.241433655332344
$ws.Range("G14").Value = 6.842156077834716
$ws.Range("H14").Value = 4.332913885022492
$ws.Range("J14").Value = 0.4094420412342998
$ws.Range("L14").Value = 0.1092732891106643
$ws.Range("N14").Value = 4.371107314139294
# Row 15
$ws.Range("C15").Value = 0.3841787595909523
$ws.Range("D15").Value = 0.03821884648418461
$ws.Range("E15").Value = 0.1392308239119444
$ws.Range("F15").Value = 7.220663700067291
$ws.Range("G15").Value = 6.817517457855274
$ws.Range("H15").Value = 4.322580703489336
$ws.Range("J15").Value = 0.4089265095695254
$ws.Range("L15").Value = 0.1093347656319708
$ws.Range("N15").Value = 4.329286057409945
# Row 16
$ws.Range("C16").Value = 0.3821542357728163
$ws.Range("D16").Value = 0.03640921834443134
$ws.Range("E16").Value = 0.1392932412751655
$ws.Range("F16").Value = 7.103146468933517
$ws.Range("G16").Value = 6.677815739129528
$ws.Range("H16").Value = 4.264236807723648
$ws.Range("J16").Value = 0.4060457876422987
$ws.Range("L16").Value = 0.1096959108538726
$ws.Range("N16").Value = 4.089429168003562
# Row 17
$ws.Range("C17").Value = 0.380969759335386
$ws.Range("D17").Value = 0.03529866632443657
$ws.Range("E17").Value = 0.1393475253000567
$ws.Range("F17").Value = 7.032385468442385
$ws.Range("G17").Value = 6.593430308115501
$ws.Range("H17").Value = 4.229215660517411
$ws.Range("J17").Value = 0.4043438101738843
$ws.Range("L17").Value = 0.109925367196098
$ws.Range("N17").Value = 3.94211849063862
# Row 18
$ws.Range("C18").Value = 0.3803094957090565
$ws.Range("D18").Value = 0.03465969297023719
$ws.Range("E18").Value = 0.1393846128102574
$ws.Range("F18").Value = 6.992170153445443
$ws.Range("G18").Value = 6.545372367118489
$ws.Range("H18").Value = 4.209353153736913
$ws.Range("J18").Value = 0.4033886989461308
$ws.Range("L18").Value = 0.1100602536844804
$ws.Range("N18").Value = 3.857331695637754
# Row 19
$ws.Range("C19").Value = 0.3800895404093296
$ws.Range("D19").Value = 0.03444330486869518
$ws.Range("E19").Value = 0.1393981761159253
$ws.Range("F19").Value = 6.978636693420697
$ws.Range("G19").Value = 6.529182394438465
$ws.Range("H19").Value = 4.202676023481445
$ws.Range("J19").Value = 0.4030693911607131
$ws.Range("L19").Value = 0.1101064239886753
$ws.Range("N19").Value = 3.828614786364199
# Row 20
$ws.Range("C20").Value = 0.3810936712141313
$ws.Range("D20").Value = 0.03541690709860035
$ws.Range("E20").Value = 0.139341139411755
$ws.Range("F20").Value = 7.039867838568171
$ws.Range("G20").Value = 6.602363649255153
$ws.Range("H20").Value = 4.232914609046873
$ws.Range("J20").Value = 0.4045225196344688
$ws.Range("L20").Value = 0.1099006401707037
$ws.Range("N20").Value = 3.95780600327754
# Row 21
$ws.Range("C21").Value = 0.3847193259934159
$ws.Range("D21").Value = 0.03868598548244506
$ws.Range("E21").Value = 0.1392196638943375
$ws.Range("F21").Value = 7.251420251305944
$ws.Range("G21").Value = 6.853997568860336
$ws.Range("H21").Value = 4.337884434673185
$ws.Range("J21").Value = 0.4096905600656129
$ws.Range("L21").Value = 0.1092439658245432
$ws.Range("N21").Value = 4.391158149571083
# Row 22
$ws.Range("C22").Value = 0.3872785531142995
$ws.Range("D22").Value = 0.04082200908542433
$ws.Range("E22").Value = 0.139192771249089
$ws.Range("F22").Value = 7.394103624924867
$ws.Range("G22").Value = 7.02283649814234
$ws.Range("H22").Value = 4.409043827508242
$ws.Range("J22").Value = 0.4132837161551492
$ws.Range("L22").Value = 0.1088406128403481
$ws.Range("N22").Value = 4.67379181795809
# Row 23
$ws.Range("C23").Value = 0.3858951678991502
$ws.Range("D23").Value = 0.03968195115069051
$ws.Range("E23").Value = 0.1392023070030604
$ws.Range("F23").Value = 7.317540962805481
$ws.Range("G23").Value = 6.93231703334925
$ws.Range("H23").Value = 4.370828136050818
$ws.Range("J23").Value = 0.4113460896759875
$ws.Range("L23").Value = 0.1090535312789793
$ws.Range("N23").Value = 4.523002190001307
# Row 24
$ws.Range("C24").Value = 0.3810375861811224
$ws.Range("D24").Value = 0.03536345202419966
$ws.Range("E24").Value = 0.1393440081633202
$ws.Range("F24").Value = 7.036483608324176
$ws.Range("G24").Value = 6.598323465582098
$ws.Range("H24").Value = 4.231241468999428
$ws.Range("J24").Value = 0.4044416523062608
$ws.Range("L24").Value = 0.1099118100100451
$ws.Range("N24").Value = 3.950713976768498
# Row 25
$ws.Range("C25").Value = 0.3765999893398657
$ws.Range("D25").Value = 0.03069630908392895
$ws.Range("E25").Value = 0.1397209344320949
$ws.Range("F25").Value = 6.751782705295085
$ws.Range("G25").Value = 6.256275961389804
$ws.Range("H25").Value = 4.091376358756406
$ws.Range("J25").Value = 0.3979029319190772
$ws.Range("L25").Value = 0.1109488849940234
$ws.Range("N25").Value = 3.331249627311138
